$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price column (D) contains text-like values (e.g. "51.633.70") that Excel would
# otherwise auto-convert to numbers; force the whole column to Text first, write
# the values, then restore the default "Normal" style so cells keep looking like
# the originals (no explicit number format), while the stored content stays text.
$priceCol = $ws.Range("D2:D51")
$priceCol.NumberFormat = "@"

$ws.Range("D2").Value = "51.633.70"
$ws.Range("E2").Value = "  +0.59%  "
$ws.Range("D3").Value = "2.996.86"
$ws.Range("E3").Value = "  +2.72%  "
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").Value = "381.38"
$ws.Range("E5").Value = "  +4.63%  "
$ws.Range("D6").Value = "106.16"
$ws.Range("E6").Value = "  +2.38%  "
$ws.Range("D7").Value = "0.546"
$ws.Range("E7").Value = "  +1.19%  "
$ws.Range("E8").Value = "  +0.07%  "
$ws.Range("E9").Value = "  +1.27%  "
$ws.Range("D10").Value = "37.59"
$ws.Range("E10").Value = "  +1.61%  "
$ws.Range("E11").Value = "  +0.60%  "
$ws.Range("D12").Value = "0.0845"
$ws.Range("E12").Value = "  +1.21%  "
$ws.Range("E13").Value = "  +1.55%  "
$ws.Range("D14").Value = "3.471.12"
$ws.Range("E14").Value = "  +2.97%  "
$ws.Range("D15").Value = "7.53"
$ws.Range("E15").Value = "  +2.11%  "
$ws.Range("D16").Value = "3.004.37"
$ws.Range("E16").Value = "  +2.96%  "
$ws.Range("D17").Value = "0.971"
$ws.Range("E17").Value = "  +1.82%  "
$ws.Range("D18").Value = "51.691.76"
$ws.Range("E18").Value = "  +0.94%  "
$ws.Range("E19").Value = "  +1.85%  "
$ws.Range("D20").Value = "7.43"
$ws.Range("E20").Value = "  +2.31%  "
$ws.Range("D21").Value = "13.09"
$ws.Range("E21").Value = "  +0.26%  "
$ws.Range("D22").Value = "0.0₃0960"
$ws.Range("E22").Value = "  +1.19%  "
$ws.Range("D23").Value = "69.40"
$ws.Range("E23").Value = "  +1.69%  "
$ws.Range("D24").Value = "264.23"
$ws.Range("E24").Value = "  +1.49%  "
$ws.Range("E25").Value = "  +3.90%  "
$ws.Range("E26").Value = "  -1.76%  "
$ws.Range("D27").Value = "7.26"
$ws.Range("E27").Value = "  +17.85%  "
$ws.Range("D28").Value = "7.49"
$ws.Range("E28").Value = "  +3.91%  "
$ws.Range("D29").Value = "26.15"
$ws.Range("E29").Value = "  +0.82%  "
$ws.Range("E30").Value = "  -0.08%  "
$ws.Range("D31").Value = "0.108"
$ws.Range("E31").Value = "  +2.48%  "
$ws.Range("E32").Value = "  -0.18%  "
$ws.Range("D33").Value = "34.81"
$ws.Range("E33").Value = "  -0.93%  "
$ws.Range("B34").Value = "VeChain"
$ws.Range("C34").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D34").Value = "0.0463"
$ws.Range("E34").Value = "  +9.22%  "
$ws.Range("B35").Value = "OKB"
$ws.Range("C35").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D35").Value = "51.46"
$ws.Range("E35").Value = "  +1.30%  "
$ws.Range("B36").Value = "Toncoin"
$ws.Range("C36").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D36").Value = "2.09"
$ws.Range("E36").Value = "  -2.60%  "
$ws.Range("E37").Value = "  +0.23%  "
$ws.Range("E38").Value = "  -0.33%  "
$ws.Range("D39").Value = "17.55"
$ws.Range("E39").Value = "  +3.25%  "
$ws.Range("E40").Value = "  -6.51%  "
$ws.Range("E41").Value = "  -0.30%  "
$ws.Range("E42").Value = "  +2.45%  "
$ws.Range("D43").Value = "123.85"
$ws.Range("E43").Value = "  +4.16%  "
$ws.Range("D44").Value = "22.25"
$ws.Range("E44").Value = "  -0.94%  "
$ws.Range("E45").Value = "  -1.47%  "
$ws.Range("E46").Value = "  +16.71%  "
$ws.Range("D47").Value = "2.062.57"
$ws.Range("E47").Value = "  +0.04%  "
$ws.Range("D48").Value = "3.27"
$ws.Range("E48").Value = "  +2.35%  "
$ws.Range("D49").Value = "2.34"
$ws.Range("E49").Value = "  +3.33%  "
$ws.Range("D50").Value = "0.0353"
$ws.Range("E50").Value = "  +12.00%  "
$ws.Range("D51").Value = "5.20"
$ws.Range("E51").Value = "  +3.27%  "

$priceCol.Style = "Normal"
